$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Convert B2, C2, D2 from static values to "today"-driven formulas ---
# (doing this first frees the now-unused "2020" shared string slot so the
# new strings we add below fill in starting at that freed index, matching
# the order Excel produced them in)
$ws.Range("B2").Formula = "=DAY(TODAY())"
$ws.Range("B2").NumberFormat = "General"

$ws.Range("C2").Formula = "=TEXT(TODAY(),""mm"")"
$ws.Range("C2").NumberFormat = "mm-dd-yy"

$ws.Range("D2").Formula = "=YEAR(TODAY())"
$ws.Range("D2").NumberFormat = "General"

# --- New "date type" columns (E:J) on row 1 - new shared strings ---
$ws.Range("F1").Value = "8 Weeks"
$ws.Range("F1").NumberFormat = "@"

$ws.Range("G1").Value = "12 Weeks"
$ws.Range("G1").NumberFormat = "@"

$ws.Range("H1").Value = "5 Weeks"
$ws.Range("H1").NumberFormat = "@"

$ws.Range("I1").Value = "6 Months Plus 8 Weeks"
$ws.Range("I1").NumberFormat = "@"

$ws.Range("J1").Value = "6 Months Plus 5 Weeks"
$ws.Range("J1").NumberFormat = "@"

$ws.Range("E1").Value = "6 Months"
$ws.Range("E1").NumberFormat = "@"

# --- New formulas on row 2 ---
$ws.Range("E2").Formula = "=DATEDIF(TODAY(),EDATE(TODAY(),6),""D"")"
$ws.Range("E2").NumberFormat = "General"

$ws.Range("F2").Formula = "=SUM(8*7)"
$ws.Range("G2").Formula = "=SUM(12*7)"
$ws.Range("H2").Formula = "=SUM(5*7)"
$ws.Range("I2").Formula = "=SUM(E2+F2)"
$ws.Range("J2").Formula = "=SUM(E2+(5*7))"

# --- Column widths (best approximation; engine snaps to 1/6-char steps) ---
$ws.Columns.Item(1).ColumnWidth = 16.666666666666668
$ws.Columns.Item(2).ColumnWidth = 13.5
$ws.Columns.Item(3).ColumnWidth = 9.666666666666666
$ws.Columns.Item(4).ColumnWidth = 4.166666666666667
$ws.Columns.Item(5).ColumnWidth = 9.666666666666666
$ws.Columns.Item(9).ColumnWidth = 18.833333333333332
$ws.Columns.Item(10).ColumnWidth = 19.0

# --- Selection moves to E1 after the edit ---
$ws.Range("E1").Select()
